# Auto-generated Excel COM-interop edit script
# Applies the horarios-141-2026-01-14 update (Linea 141 schedule refresh)
$wb = $excel.ActiveWorkbook

# --- LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:20:49"
$ws.Cells.Item(3, 1).Value = "Total filas: 253"
$ws.Cells.Item(115, 1).Value = "10:13:53"
$ws.Cells.Item(115, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(115, 4).Value = 9
$ws.Cells.Item(116, 1).Value = "08:31:53"
$ws.Cells.Item(116, 3).Value = "17_ROMERO"
$ws.Cells.Item(116, 4).Value = 111
$ws.Cells.Item(162, 1).Value = "10:13:53"
$ws.Cells.Item(162, 3).Value = "14_ABASTO"
$ws.Cells.Item(162, 4).Value = 113
$ws.Cells.Item(163, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(164, 1).Value = "10:52:37"
$ws.Cells.Item(164, 3).Value = "10_OLMOS"
$ws.Cells.Item(164, 4).Value = 74
$ws.Cells.Item(178, 1).Value = "11:46:46"
$ws.Cells.Item(178, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(178, 4).Value = 48
$ws.Cells.Item(179, 1).Value = "11:17:39"
$ws.Cells.Item(179, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(179, 4).Value = 77
$ws.Cells.Item(193, 1).Value = "10:52:37"
$ws.Cells.Item(193, 3).Value = "15_ABASTO"
$ws.Cells.Item(193, 4).Value = 118
$ws.Cells.Item(194, 1).Value = "12:50:41"
$ws.Cells.Item(194, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(195, 1).Value = "12:35:30"
$ws.Cells.Item(195, 3).Value = "14_ABASTO"
$ws.Cells.Item(195, 4).Value = 27
$ws.Cells.Item(196, 1).Value = "12:01:11"
$ws.Cells.Item(196, 3).Value = "15_ABASTO"
$ws.Cells.Item(196, 4).Value = 61
$ws.Cells.Item(233, 1).Value = "14:20:49"
$ws.Cells.Item(233, 2).Value = "14:28"
$ws.Cells.Item(233, 3).Value = "15_ABASTO"
$ws.Cells.Item(233, 4).Value = 8
$ws.Cells.Item(234, 1).Value = "13:51:32"
$ws.Cells.Item(234, 2).Value = "14:30"
$ws.Cells.Item(234, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(234, 4).Value = 39
$ws.Cells.Item(235, 1).Value = "12:50:41"
$ws.Cells.Item(235, 2).Value = "14:44"
$ws.Cells.Item(235, 3).Value = "14_ABASTO"
$ws.Cells.Item(235, 4).Value = 114
$ws.Cells.Item(236, 1).Value = "14:20:49"
$ws.Cells.Item(236, 2).Value = "14:46"
$ws.Cells.Item(236, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(236, 4).Value = 26
$ws.Cells.Item(237, 2).Value = "14:56"
$ws.Cells.Item(237, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(237, 4).Value = 98
$ws.Cells.Item(238, 2).Value = "14:58"
$ws.Cells.Item(238, 3).Value = "215B_EL PATO"
$ws.Cells.Item(238, 4).Value = 100
$ws.Cells.Item(239, 1).Value = "13:18:40"
$ws.Cells.Item(239, 2).Value = "15:00"
$ws.Cells.Item(239, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(239, 4).Value = 102
$ws.Cells.Item(240, 2).Value = "15:05"
$ws.Cells.Item(240, 3).Value = "10_OLMOS"
$ws.Cells.Item(240, 4).Value = 107
$ws.Cells.Item(241, 2).Value = "15:10"
$ws.Cells.Item(241, 3).Value = "17_ROMERO"
$ws.Cells.Item(241, 4).Value = 79
$ws.Cells.Item(242, 1).Value = "13:18:40"
$ws.Cells.Item(242, 2).Value = "15:13"
$ws.Cells.Item(242, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(242, 4).Value = 115
$ws.Cells.Item(243, 2).Value = "15:20"
$ws.Cells.Item(243, 3).Value = "15_ABASTO"
$ws.Cells.Item(243, 4).Value = 89
$ws.Cells.Item(244, 1).Value = "14:20:49"
$ws.Cells.Item(244, 2).Value = "15:21"
$ws.Cells.Item(244, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(244, 4).Value = 61
$ws.Cells.Item(245, 2).Value = "15:26"
$ws.Cells.Item(245, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(245, 4).Value = 95
$ws.Cells.Item(246, 2).Value = "15:32"
$ws.Cells.Item(246, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(246, 4).Value = 101
$ws.Cells.Item(247, 2).Value = "15:34"
$ws.Cells.Item(247, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(247, 4).Value = 103
$ws.Cells.Item(248, 1).Value = "14:20:49"
$ws.Cells.Item(248, 2).Value = "15:36"
$ws.Cells.Item(248, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(248, 4).Value = 76
$ws.Cells.Item(249, 1).Value = "13:51:32"
$ws.Cells.Item(249, 2).Value = "15:37"
$ws.Cells.Item(249, 3).Value = "10_OLMOS"
$ws.Cells.Item(249, 4).Value = 106
$ws.Cells.Item(249, 5).Value = "LP1912"
$ws.Cells.Item(250, 1).Value = "13:51:32"
$ws.Cells.Item(250, 2).Value = "15:38"
$ws.Cells.Item(250, 3).Value = "215A_EL PATO"
$ws.Cells.Item(250, 4).Value = 107
$ws.Cells.Item(250, 5).Value = "LP1912"
$ws.Cells.Item(251, 1).Value = "14:20:49"
$ws.Cells.Item(251, 2).Value = "15:45"
$ws.Cells.Item(251, 3).Value = "14_ABASTO"
$ws.Cells.Item(251, 4).Value = 85
$ws.Cells.Item(251, 5).Value = "LP1912"
$ws.Cells.Item(252, 1).Value = "13:51:32"
$ws.Cells.Item(252, 2).Value = "15:46"
$ws.Cells.Item(252, 3).Value = "14_ABASTO"
$ws.Cells.Item(252, 4).Value = 115
$ws.Cells.Item(252, 5).Value = "LP1912"
$ws.Cells.Item(253, 1).Value = "13:51:32"
$ws.Cells.Item(253, 2).Value = "15:46"
$ws.Cells.Item(253, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(253, 4).Value = 115
$ws.Cells.Item(253, 5).Value = "LP1912"
$ws.Cells.Item(254, 1).Value = "14:20:49"
$ws.Cells.Item(254, 2).Value = "15:53"
$ws.Cells.Item(254, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(254, 4).Value = 93
$ws.Cells.Item(254, 5).Value = "LP1912"
$ws.Cells.Item(255, 1).Value = "14:20:49"
$ws.Cells.Item(255, 2).Value = "15:55"
$ws.Cells.Item(255, 3).Value = "17_ROMERO"
$ws.Cells.Item(255, 4).Value = 95
$ws.Cells.Item(255, 5).Value = "LP1912"
$ws.Cells.Item(256, 1).Value = "14:20:49"
$ws.Cells.Item(256, 2).Value = "15:56"
$ws.Cells.Item(256, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(256, 4).Value = 96
$ws.Cells.Item(256, 5).Value = "LP1912"
$ws.Cells.Item(257, 1).Value = "14:20:49"
$ws.Cells.Item(257, 2).Value = "16:04"
$ws.Cells.Item(257, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(257, 4).Value = 104
$ws.Cells.Item(257, 5).Value = "LP1912"
$ws.Cells.Item(258, 1).Value = "14:20:49"
$ws.Cells.Item(258, 2).Value = "16:14"
$ws.Cells.Item(258, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(258, 4).Value = 114
$ws.Cells.Item(258, 5).Value = "LP1912"

# --- LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:20:49"

# --- 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:20:49"
$ws.Cells.Item(3, 1).Value = "Total filas: 38"
$ws.Cells.Item(43, 1).Value = "14:20:49"
$ws.Cells.Item(43, 2).Value = "16:13"
$ws.Cells.Item(43, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(43, 4).Value = 113
$ws.Cells.Item(43, 5).Value = "L6203"
